$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "error solve ifrs list" -- refresh IFRS-consolidated financial figures
# for rows 2-9 (FY2014 .. FY2021E) with corrected source data.

# Row 2
$ws.Range("D2").Value = 16347
$ws.Range("E2").Value = 1527
$ws.Range("F2").Value = 1723
$ws.Range("G2").Value = 1310
$ws.Range("H2").Value = 1163
$ws.Range("I2").Value = 1095
$ws.Range("J2").Value = 68
$ws.Range("K2").Value = 14601
$ws.Range("L2").Value = 5443
$ws.Range("M2").Value = 9157
$ws.Range("N2").Value = 8734
$ws.Range("O2").Value = 423
$ws.Range("P2").Value = 190
$ws.Range("Q2").Value = 2817
$ws.Range("R2").Value = -4546
$ws.Range("S2").Value = 1523
$ws.Range("T2").Value = 1431
$ws.Range("U2").Value = 1386
$ws.Range("V2").Value = 1900
$ws.Range("W2").Value = 9.34
$ws.Range("X2").Value = 7.12
$ws.Range("Y2").Value = 13.04
$ws.Range("Z2").Value = 8.9
$ws.Range("AA2").Value = 59.44
$ws.Range("AB2").Value = 5299.62
$ws.Range("AC2").Value = 2882
$ws.Range("AD2").Value = 24.71
$ws.Range("AE2").Value = 25834
$ws.Range("AF2").Value = 2.76
$ws.Range("AG2").Value = 1150
$ws.Range("AH2").Value = 1.62
$ws.Range("AI2").Value = 35.51
$ws.Range("AJ2").Value = 37999178

# Row 3
$ws.Range("D3").Value = 17996
$ws.Range("E3").Value = 1726
$ws.Range("F3").Value = 1733
$ws.Range("G3").Value = 1562
$ws.Range("H3").Value = 1554
$ws.Range("I3").Value = 1536
$ws.Range("J3").Value = 18
$ws.Range("K3").Value = 13944
$ws.Range("L3").Value = 4141
$ws.Range("M3").Value = 9804
$ws.Range("N3").Value = 9801
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 190
$ws.Range("Q3").Value = 3179
$ws.Range("R3").Value = -324
$ws.Range("S3").Value = -2297
$ws.Range("T3").Value = 1371
$ws.Range("U3").Value = 1808
$ws.Range("W3").Value = 9.59
$ws.Range("X3").Value = 8.640000000000001
$ws.Range("Y3").Value = 16.58
$ws.Range("Z3").Value = 10.89
$ws.Range("AA3").Value = 42.24
$ws.Range("AB3").Value = 5862.26
$ws.Range("AC3").Value = 4043
$ws.Range("AD3").Value = 24.66
$ws.Range("AE3").Value = 28988
$ws.Range("AF3").Value = 3.44
$ws.Range("AG3").Value = 1200
$ws.Range("AH3").Value = 1.2
$ws.Range("AI3").Value = 26.41
$ws.Range("AJ3").Value = 37999178

# Row 4
$ws.Range("D4").Value = 18302
$ws.Range("E4").Value = 2057
$ws.Range("F4").Value = 2057
$ws.Range("G4").Value = 1862
$ws.Range("H4").Value = 1405
$ws.Range("I4").Value = 1405
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 14604
$ws.Range("L4").Value = 3604
$ws.Range("M4").Value = 10999
$ws.Range("N4").Value = 10999
$ws.Range("P4").Value = 190
$ws.Range("Q4").Value = 2732
$ws.Range("R4").Value = -2135
$ws.Range("S4").Value = -410
$ws.Range("T4").Value = 1336
$ws.Range("U4").Value = 1396
$ws.Range("W4").Value = 11.24
$ws.Range("X4").Value = 7.68
$ws.Range("Y4").Value = 13.51
$ws.Range("Z4").Value = 9.85
$ws.Range("AA4").Value = 32.77
$ws.Range("AB4").Value = 6491.37
$ws.Range("AC4").Value = 3698
$ws.Range("AD4").Value = 23.71
$ws.Range("AE4").Value = 32533
$ws.Range("AF4").Value = 2.7
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 1.43
$ws.Range("AI4").Value = 30.07
$ws.Range("AJ4").Value = 37999178

# Row 5
$ws.Range("D5").Value = 19423
$ws.Range("E5").Value = 2026
$ws.Range("F5").Value = 2026
$ws.Range("G5").Value = 1914
$ws.Range("H5").Value = 1433
$ws.Range("I5").Value = 1433
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 15815
$ws.Range("L5").Value = 3737
$ws.Range("M5").Value = 12078
$ws.Range("N5").Value = 12078
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 190
$ws.Range("Q5").Value = 2942
$ws.Range("R5").Value = -2247
$ws.Range("S5").Value = -423
$ws.Range("T5").Value = 1221
$ws.Range("U5").Value = 1721
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 10.43
$ws.Range("X5").Value = 7.38
$ws.Range("Y5").Value = 12.42
$ws.Range("Z5").Value = 9.42
$ws.Range("AA5").Value = 30.94
$ws.Range("AB5").Value = 7060.34
$ws.Range("AC5").Value = 3772
$ws.Range("AD5").Value = 28.37
$ws.Range("AE5").Value = 35724
$ws.Range("AF5").Value = 3
$ws.Range("AG5").Value = 2500
$ws.Range("AH5").Value = 2.34
$ws.Range("AI5").Value = 58.97
$ws.Range("AJ5").Value = 37999178

# Row 6
$ws.Range("D6").Value = 20183
$ws.Range("E6").Value = 1991
$ws.Range("F6").Value = 1991
$ws.Range("G6").Value = 1986
$ws.Range("H6").Value = 1030
$ws.Range("I6").Value = 1030
$ws.Range("K6").Value = 16829
$ws.Range("L6").Value = 4205
$ws.Range("M6").Value = 12623
$ws.Range("N6").Value = 12623
$ws.Range("P6").Value = 190
$ws.Range("Q6").Value = 2370
$ws.Range("R6").Value = -1762
$ws.Range("S6").Value = -845
$ws.Range("T6").Value = 1300
$ws.Range("U6").Value = 1070
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 9.869999999999999
$ws.Range("X6").Value = 5.1
$ws.Range("Y6").Value = 8.34
$ws.Range("Z6").Value = 6.31
$ws.Range("AA6").Value = 33.31
$ws.Range("AB6").Value = 7347.08
$ws.Range("AC6").Value = 2710
$ws.Range("AD6").Value = 37.09
$ws.Range("AE6").Value = 37335
$ws.Range("AF6").Value = 2.69
$ws.Range("AG6").Value = 2500
$ws.Range("AH6").Value = 2.49
$ws.Range("AI6").Value = 82.09
$ws.Range("AJ6").Value = 37999178

# Row 7
$ws.Range("D7").Value = 21337
$ws.Range("E7").Value = 2068
$ws.Range("G7").Value = 2037
$ws.Range("H7").Value = 1550
$ws.Range("I7").Value = 1549
$ws.Range("K7").Value = 17646
$ws.Range("L7").Value = 4501
$ws.Range("M7").Value = 13146
$ws.Range("N7").Value = 13146
$ws.Range("P7").Value = 190
$ws.Range("Q7").Value = 2740
$ws.Range("R7").Value = -1288
$ws.Range("S7").Value = -890
$ws.Range("T7").Value = 1387
$ws.Range("U7").Value = 1366
$ws.Range("W7").Value = 9.69
$ws.Range("X7").Value = 7.27
$ws.Range("Y7").Value = 12.02
$ws.Range("Z7").Value = 8.99
$ws.Range("AA7").Value = 34.24
$ws.Range("AC7").Value = 4076
$ws.Range("AD7").Value = 22.72
$ws.Range("AE7").Value = 38881
$ws.Range("AF7").Value = 2.38
$ws.Range("AG7").Value = 2600
$ws.Range("AH7").Value = 2.81
$ws.Range("AI7").Value = 63.78

# Row 8
$ws.Range("D8").Value = 22607
$ws.Range("E8").Value = 2235
$ws.Range("G8").Value = 2215
$ws.Range("H8").Value = 1685
$ws.Range("I8").Value = 1685
$ws.Range("K8").Value = 18220
$ws.Range("L8").Value = 4710
$ws.Range("M8").Value = 13510
$ws.Range("N8").Value = 13510
$ws.Range("P8").Value = 190
$ws.Range("Q8").Value = 2915
$ws.Range("R8").Value = -920
$ws.Range("S8").Value = -870
$ws.Range("T8").Value = 1400
$ws.Range("U8").Value = 1940
$ws.Range("W8").Value = 9.890000000000001
$ws.Range("X8").Value = 7.45
$ws.Range("Y8").Value = 12.77
$ws.Range("Z8").Value = 9.48
$ws.Range("AA8").Value = 34.86
$ws.Range("AC8").Value = 4434
$ws.Range("AD8").Value = 20.54
$ws.Range("AE8").Value = 39958
$ws.Range("AF8").Value = 2.28
$ws.Range("AG8").Value = 2650
$ws.Range("AH8").Value = 2.91
$ws.Range("AI8").Value = 59.76

# Row 9
$ws.Range("D9").Value = 23620
$ws.Range("E9").Value = 2395
$ws.Range("G9").Value = 2385
$ws.Range("H9").Value = 1810
$ws.Range("I9").Value = 1810
$ws.Range("K9").Value = 19150
$ws.Range("L9").Value = 4905
$ws.Range("M9").Value = 14245
$ws.Range("N9").Value = 14245
$ws.Range("P9").Value = 190
$ws.Range("Q9").Value = 3250
$ws.Range("R9").Value = -820
$ws.Range("S9").Value = -905
$ws.Range("T9").Value = 1400
$ws.Range("U9").Value = 2040
$ws.Range("W9").Value = 10.14
$ws.Range("X9").Value = 7.66
$ws.Range("Y9").Value = 13.04
$ws.Range("Z9").Value = 9.69
$ws.Range("AA9").Value = 34.43
$ws.Range("AC9").Value = 4763
$ws.Range("AD9").Value = 19.13
$ws.Range("AE9").Value = 42132
$ws.Range("AF9").Value = 2.16
$ws.Range("AG9").Value = 2750
$ws.Range("AH9").Value = 3.02
$ws.Range("AI9").Value = 57.73

# These cells have no corresponding data in the corrected source and must
# be cleared entirely (not just zeroed) so the cell element is dropped.
$ws.Range("V3").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()
